$wb = $excel.ActiveWorkbook

# Sheet 1: VENTAS POR GRUPO
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M13").Value = 1221.56
$ws1.Range("D20").Value = 183.17
$ws1.Range("D30").Value = "2 de 28"

# Sheet 2: VENTA MENSUAL
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F13").Value = 1227.32
$ws2.Range("F20").Value = 765.7
$ws2.Range("F30").Value = 6519.4

# Sheet 3: CUMPLIMIENTO MENSUAL
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D3").Value = 1165.63
$ws3.Range("E3").Value = 1954.4845
$ws3.Range("F3").Value = 0.3735856488600018

$ws3.Range("D16").Value = 3961.67
$ws3.Range("E16").Value = 14836.94
$ws3.Range("F16").Value = 0.2107427091683906

$ws3.Range("D19").Value = 6513.64
$ws3.Range("E19").Value = 23024.15107555787
$ws3.Range("F19").Value = 0.2205188594955548
